# Weekly data refresh: insert this week's new Nectarín price rows
# (Macroferia Regional de Talca, date 2022-12-23 / serial 44918) at the
# top of the data block (row 685), pushing the existing historical rows
# down by 4 (they keep all their original values, just shifted to
# rows 689-764).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 685, shifting everything
# from 685 downward (previously ending at 760) down to end at 764.
$ws.Rows("685:688").Insert()

# Columns that are constant for every data row in this sheet (it's a
# single market / region / product subset export).
$constA = 5
$constB = "Macroferia Regional de Talca"
$constC = "Maule"
$constE = 7
$constF = "Fruta"
$constG = 100103
$constH = "Frutos de hueso (carozo)"
$constI = 100103006
$constJ = "Nectarín"
$constR = "Región de O'Higgins"

# New-row data: Date, Variedad (K), Calidad (L), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P), Unidad (Q), Precio $/Kg (S), Kg/unidad (T)
$newRows = @(
    @{ D = 44918; K = "Artic Star";  L = "Especial"; M = 240; N = 14000; O = 14000; P = 14000; Q = "$/bandeja 15 kilos granel";   S = 933; T = 15 },
    @{ D = 44918; K = "Artic Star";  L = "Primera";  M = 200; N = 12000; O = 12000; P = 12000; Q = "$/bandeja 15 kilos granel";   S = 800; T = 15 },
    @{ D = 44918; K = "Super Queen"; L = "Especial"; M = 600; N = 12000; O = 13000; P = 12583; Q = "$/bandeja 15 kilos empedrada"; S = 839; T = 15 },
    @{ D = 44918; K = "Super Queen"; L = "Primera";  M = 180; N = 11000; O = 11000; P = 11000; Q = "$/bandeja 15 kilos empedrada"; S = 733; T = 15 }
)

$r = 685
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2  = $constA
    $ws.Cells.Item($r, 2).Value2  = $constB
    $ws.Cells.Item($r, 3).Value2  = $constC
    $ws.Cells.Item($r, 4).Value2  = $row.D
    $ws.Cells.Item($r, 5).Value2  = $constE
    $ws.Cells.Item($r, 6).Value2  = $constF
    $ws.Cells.Item($r, 7).Value2  = $constG
    $ws.Cells.Item($r, 8).Value2  = $constH
    $ws.Cells.Item($r, 9).Value2  = $constI
    $ws.Cells.Item($r, 10).Value2 = $constJ
    $ws.Cells.Item($r, 11).Value2 = $row.K
    $ws.Cells.Item($r, 12).Value2 = $row.L
    $ws.Cells.Item($r, 13).Value2 = $row.M
    $ws.Cells.Item($r, 14).Value2 = $row.N
    $ws.Cells.Item($r, 15).Value2 = $row.O
    $ws.Cells.Item($r, 16).Value2 = $row.P
    $ws.Cells.Item($r, 17).Value2 = $row.Q
    $ws.Cells.Item($r, 18).Value2 = $constR
    $ws.Cells.Item($r, 19).Value2 = $row.S
    $ws.Cells.Item($r, 20).Value2 = $row.T
    $r++
}
